# Apply crypto price/volume table refresh (GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.890.46"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.892.97"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7708"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.07"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3123"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.66"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07218"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08700"
$ws.Range("E11").Value = "  +7.52%  "
$ws.Range("D12").Value = "2.106.20"
$ws.Range("E12").Value = "  +10.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7704"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.20"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.212"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "30.158.95"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.92"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.92"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.291.68"
$ws.Range("E20").Value = "  +4.83%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007863"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.161"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1597"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.512"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.23"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.036"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.434"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.526"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.117"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05467"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.247"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7561"
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.60"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "1.095.71"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.031"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8556"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.97"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.883"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.161.21"
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.619"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.830"
$ws.Range("E51").Value = "  -2.38%  "
